$wb = $excel.ActiveWorkbook

# --- Arkusz1 (sheet1): clear the now-unused helper/legend cell text, ---
# --- keeping styles intact, and move the selection to H21.           ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1").Value = $null
$ws1.Range("A2").Value = $null
$ws1.Range("A3").Value = $null

$ws1.Range("F10").Value = $null
$ws1.Range("F11:J11").Value = $null
$ws1.Range("F12:J12").Value = $null
$ws1.Range("F13:H13").Value = $null
$ws1.Range("J13").Value = $null
$ws1.Range("F14:H14").Value = $null
$ws1.Range("J14").Value = $null

# Row 12 no longer wraps multi-line text once its cells are cleared, so
# it shrinks back down to the same (thick-bottom-border) row height used
# by its neighbouring rows.
$ws1.Rows.Item(12).RowHeight = 18.600000000000001

# Move the selection on this sheet (this also temporarily activates the
# sheet, which we correct below by activating Arkusz4 last).
$ws1.Range("H21").Select()

# --- Active tab moves from Arkusz3 to Arkusz4. ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
